$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = '中芯国际'
$ws.Range('A3').Value = '浪潮信息'
$ws.Range('C3').Value = '紫金矿业'
$ws.Range('A4').Value = '三花智控'
$ws.Range('B4').Value = '中芯国际'
$ws.Range('C4').Value = '上海电气'
$ws.Range('A5').Value = '上海电气'
$ws.Range('B5').Value = '江西铜业'
$ws.Range('C5').Value = '深科技'
$ws.Range('A6').Value = '深科技'
$ws.Range('B6').Value = '特变电工'
$ws.Range('C6').Value = '中芯国际'
$ws.Range('A7').Value = '山子高科'
$ws.Range('B7').Value = '东方财富'
$ws.Range('C7').Value = '大众公用'
$ws.Range('A8').Value = 'XD紫金矿'
$ws.Range('B8').Value = '山子高科'
$ws.Range('C8').Value = '三花智控'
$ws.Range('A9').Value = '江西铜业'
$ws.Range('B9').Value = '浪潮信息'
$ws.Range('C9').Value = '赣锋锂业'
$ws.Range('A10').Value = '永鼎股份'
$ws.Range('B10').Value = '宁德时代'
$ws.Range('C10').Value = '中电鑫龙'
$ws.Range('A11').Value = '大众公用'
$ws.Range('B11').Value = 'XD紫金矿'
$ws.Range('C11').Value = '江西铜业'
$ws.Range('A12').Value = '兆易创新'
$ws.Range('B12').Value = '合锻智能'
$ws.Range('C12').Value = '浪潮信息'
$ws.Range('A13').Value = '合锻智能'
$ws.Range('B13').Value = '贵州茅台'
$ws.Range('C13').Value = '洛阳钼业'
$ws.Range('A14').Value = '赣锋锂业'
$ws.Range('B14').Value = '中国瑞林'
$ws.Range('C14').Value = '合锻智能'
$ws.Range('B15').Value = '华虹公司'
$ws.Range('C15').Value = '永鼎股份'
$ws.Range('A16').Value = '江波龙'
$ws.Range('B16').Value = '赣锋锂业'
$ws.Range('C16').Value = '天赐材料'
$ws.Range('A17').Value = '张江高科'
$ws.Range('B17').Value = '深科技'
$ws.Range('C17').Value = '华虹公司'
$ws.Range('A18').Value = '华虹公司'
$ws.Range('B18').Value = '大众公用'
$ws.Range('C18').Value = '张江高科'
$ws.Range('A19').Value = '洛阳钼业'
$ws.Range('B19').Value = '三花智控'
$ws.Range('C19').Value = '赛力斯'
$ws.Range('A20').Value = '东方财富'
$ws.Range('C20').Value = '华友钴业'
$ws.Range('A21').Value = '宁德时代'
$ws.Range('B21').Value = '华友钴业'
$ws.Range('C21').Value = '万向钱潮'
